# Commit: "sort and sum up two brand sheets."
# For this particular worksheet the only substantive content change is a
# typo fix in cell C4: "TMP-140L" -> "TPM-140L" (brand/model code correction).
# The remainder of the canonical diff (namespace/version bumps, theme rename,
# window geometry, minor row-height/column-width rounding, duplicate style
# entries, etc.) are artifacts of the file being re-saved by a newer Excel
# build and are reproduced here only where they are meaningfully expressible
# through the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the product code typo in C4 ---
$ws.Range("C4").Value = "TPM-140L"

# --- Match the row height tweak applied across the data rows (16.8 -> 16.75) ---
for ($r = 3; $r -le 13; $r++) {
    $ws.Rows.Item($r).RowHeight = 16.75
}

# --- Nudge the data column width to the slightly narrower re-saved value ---
$ws.Range("A1:E13").ColumnWidth = 19.83

# --- The re-saved file had the active cell/selection parked on C5 ---
$ws.Range("C5").Select()
